$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value = 0.6840000000000001
$ws.Cells.Item(2,3).Value = 0.516
$ws.Cells.Item(2,4).Value = 0.357
$ws.Cells.Item(2,5).Value = 0.093
$ws.Cells.Item(2,6).Value = 0.474
$ws.Cells.Item(2,7).Value = 0.101
$ws.Cells.Item(2,8).Value = 0.167
$ws.Cells.Item(2,9).Value = 0.699
$ws.Cells.Item(2,10).Value = 0.949
$ws.Cells.Item(2,11).Value = 0.805
$ws.Cells.Item(3,2).Value = 0.418
$ws.Cells.Item(3,4).Value = 0.357
$ws.Cells.Item(3,5).Value = 0.062
$ws.Cells.Item(3,6).Value = 0.327
$ws.Cells.Item(3,7).Value = 0.82
$ws.Cells.Item(3,8).Value = 0.468
$ws.Cells.Item(3,9).Value = 0.742
$ws.Cells.Item(3,10).Value = 0.235
$ws.Cells.Item(3,11).Value = 0.357
$ws.Cells.Item(4,2).Value = 0.446
$ws.Cells.Item(4,3).Value = 0.584
$ws.Cells.Item(4,4).Value = 0.44
$ws.Cells.Item(4,5).Value = 0.047
$ws.Cells.Item(4,6).Value = 0.327
$ws.Cells.Item(4,7).Value = 0.73
$ws.Cells.Item(4,8).Value = 0.451
$ws.Cells.Item(4,9).Value = 0.721
$ws.Cells.Item(4,10).Value = 0.316
$ws.Cells.Item(4,11).Value = 0.44
$ws.Cells.Item(5,2).Value = 0.418
$ws.Cells.Item(5,3).Value = 0.575
$ws.Cells.Item(5,4).Value = 0.366
$ws.Cells.Item(5,5).Value = 0.047
$ws.Cells.Item(5,6).Value = 0.327
$ws.Cells.Item(5,7).Value = 0.73
$ws.Cells.Item(5,8).Value = 0.451
$ws.Cells.Item(5,9).Value = 0.721
$ws.Cells.Item(5,10).Value = 0.316
$ws.Cells.Item(5,11).Value = 0.44
$ws.Cells.Item(6,2).Value = 0.695
$ws.Cells.Item(6,3).Value = 0.4
$ws.Cells.Item(6,4).Value = 0.8159999999999999
$ws.Cells.Item(6,5).Value = 0.115
$ws.Cells.Item(6,6).Value = 0.625
$ws.Cells.Item(6,7).Value = 0.056
$ws.Cells.Item(6,8).Value = 0.103
$ws.Cells.Item(6,9).Value = 0.697
$ws.Cells.Item(6,10).Value = 0.985
$ws.Cells.Item(6,11).Value = 0.8159999999999999
$ws.Cells.Item(7,2).Value = 0.453
$ws.Cells.Item(7,3).Value = 0.587
$ws.Cells.Item(7,4).Value = 0.431
$ws.Cells.Item(7,5).Value = 0.091
$ws.Cells.Item(7,6).Value = 0.338
$ws.Cells.Item(7,7).Value = 0.787
$ws.Cells.Item(7,8).Value = 0.473
$ws.Cells.Item(7,9).Value = 0.756
$ws.Cells.Item(7,10).Value = 0.301
$ws.Cells.Item(7,11).Value = 0.431
$ws.Cells.Item(8,2).Value = 0.463
$ws.Cells.Item(8,3).Value = 0.547
$ws.Cells.Item(8,4).Value = 0.488
$ws.Cells.Item(8,5).Value = 0.034
$ws.Cells.Item(8,6).Value = 0.324
$ws.Cells.Item(8,7).Value = 0.663
$ws.Cells.Item(8,8).Value = 0.435
$ws.Cells.Item(8,9).Value = 0.709
$ws.Cells.Item(8,10).Value = 0.372
$ws.Cells.Item(8,11).Value = 0.488
$ws.Cells.Item(9,2).Value = 0.453
$ws.Cells.Item(9,3).Value = 0.585
$ws.Cells.Item(9,4).Value = 0.443
$ws.Cells.Item(9,5).Value = 0.07000000000000001
$ws.Cells.Item(9,6).Value = 0.333
$ws.Cells.Item(9,7).Value = 0.753
$ws.Cells.Item(9,8).Value = 0.462
$ws.Cells.Item(9,9).Value = 0.738
$ws.Cells.Item(9,10).Value = 0.316
$ws.Cells.Item(9,11).Value = 0.443
$ws.Cells.Item(10,2).Value = 0.695
$ws.Cells.Item(10,3).Value = 0.4
$ws.Cells.Item(10,4).Value = 0.8159999999999999
$ws.Cells.Item(10,5).Value = 0.115
$ws.Cells.Item(10,6).Value = 0.625
$ws.Cells.Item(10,7).Value = 0.056
$ws.Cells.Item(10,8).Value = 0.103
$ws.Cells.Item(10,9).Value = 0.697
$ws.Cells.Item(10,10).Value = 0.985
$ws.Cells.Item(10,11).Value = 0.8159999999999999
$ws.Cells.Item(11,2).Value = 0.449
$ws.Cells.Item(11,3).Value = 0.587
$ws.Cells.Item(11,4).Value = 0.429
$ws.Cells.Item(11,5).Value = 0.079
$ws.Cells.Item(11,6).Value = 0.335
$ws.Cells.Item(11,7).Value = 0.775
$ws.Cells.Item(11,8).Value = 0.468
$ws.Cells.Item(11,9).Value = 0.747
$ws.Cells.Item(11,10).Value = 0.301
$ws.Cells.Item(11,11).Value = 0.429
$ws.Cells.Item(12,2).Value = 0.463
$ws.Cells.Item(12,3).Value = 0.547
$ws.Cells.Item(12,4).Value = 0.488
$ws.Cells.Item(12,5).Value = 0.034
$ws.Cells.Item(12,6).Value = 0.324
$ws.Cells.Item(12,7).Value = 0.663
$ws.Cells.Item(12,8).Value = 0.435
$ws.Cells.Item(12,9).Value = 0.709
$ws.Cells.Item(12,10).Value = 0.372
$ws.Cells.Item(12,11).Value = 0.488
$ws.Cells.Item(13,2).Value = 0.453
$ws.Cells.Item(13,3).Value = 0.585
$ws.Cells.Item(13,4).Value = 0.443
$ws.Cells.Item(13,5).Value = 0.07000000000000001
$ws.Cells.Item(13,6).Value = 0.333
$ws.Cells.Item(13,7).Value = 0.753
$ws.Cells.Item(13,8).Value = 0.462
$ws.Cells.Item(13,9).Value = 0.738
$ws.Cells.Item(13,10).Value = 0.316
$ws.Cells.Item(13,11).Value = 0.443
$ws.Cells.Item(14,2).Value = 0.695
$ws.Cells.Item(14,3).Value = 0.4
$ws.Cells.Item(14,4).Value = 0.8159999999999999
$ws.Cells.Item(14,5).Value = 0.115
$ws.Cells.Item(14,6).Value = 0.625
$ws.Cells.Item(14,7).Value = 0.056
$ws.Cells.Item(14,8).Value = 0.103
$ws.Cells.Item(14,9).Value = 0.697
$ws.Cells.Item(14,10).Value = 0.985
$ws.Cells.Item(14,11).Value = 0.8159999999999999
$ws.Cells.Item(15,2).Value = 0.453
$ws.Cells.Item(15,3).Value = 0.587
$ws.Cells.Item(15,4).Value = 0.431
$ws.Cells.Item(15,5).Value = 0.091
$ws.Cells.Item(15,6).Value = 0.338
$ws.Cells.Item(15,7).Value = 0.787
$ws.Cells.Item(15,8).Value = 0.473
$ws.Cells.Item(15,9).Value = 0.756
$ws.Cells.Item(15,10).Value = 0.301
$ws.Cells.Item(15,11).Value = 0.431
$ws.Cells.Item(16,2).Value = 0.463
$ws.Cells.Item(16,3).Value = 0.547
$ws.Cells.Item(16,4).Value = 0.488
$ws.Cells.Item(16,5).Value = 0.034
$ws.Cells.Item(16,6).Value = 0.324
$ws.Cells.Item(16,7).Value = 0.663
$ws.Cells.Item(16,8).Value = 0.435
$ws.Cells.Item(16,9).Value = 0.709
$ws.Cells.Item(16,10).Value = 0.372
$ws.Cells.Item(16,11).Value = 0.488
$ws.Cells.Item(17,2).Value = 0.453
$ws.Cells.Item(17,3).Value = 0.585
$ws.Cells.Item(17,4).Value = 0.443
$ws.Cells.Item(17,5).Value = 0.07000000000000001
$ws.Cells.Item(17,6).Value = 0.333
$ws.Cells.Item(17,7).Value = 0.753
$ws.Cells.Item(17,8).Value = 0.462
$ws.Cells.Item(17,9).Value = 0.738
$ws.Cells.Item(17,10).Value = 0.316
$ws.Cells.Item(17,11).Value = 0.443
$ws.Cells.Item(18,2).Value = 0.681
$ws.Cells.Item(18,4).Value = 0.788
$ws.Cells.Item(18,5).Value = 0.172
$ws.Cells.Item(18,6).Value = 0.481
$ws.Cells.Item(18,7).Value = 0.281
$ws.Cells.Item(18,8).Value = 0.355
$ws.Cells.Item(18,9).Value = 0.725
$ws.Cells.Item(18,10).Value = 0.862
$ws.Cells.Item(18,11).Value = 0.788
$ws.Cells.Item(19,2).Value = 0.709
$ws.Cells.Item(19,3).Value = 0.681
$ws.Cells.Item(19,4).Value = 0.801
$ws.Cells.Item(19,5).Value = 0.272
$ws.Cells.Item(19,6).Value = 0.547
$ws.Cells.Item(19,7).Value = 0.393
$ws.Cells.Item(19,8).Value = 0.458
$ws.Cells.Item(19,9).Value = 0.756
$ws.Cells.Item(19,10).Value = 0.852
$ws.Cells.Item(19,11).Value = 0.801
$ws.Cells.Item(20,2).Value = 0.705
$ws.Cells.Item(20,3).Value = 0.6840000000000001
$ws.Cells.Item(20,4).Value = 0.8
$ws.Cells.Item(20,5).Value = 0.258
$ws.Cells.Item(20,6).Value = 0.541
$ws.Cells.Item(20,7).Value = 0.371
$ws.Cells.Item(20,8).Value = 0.44
$ws.Cells.Item(20,9).Value = 0.75
$ws.Cells.Item(20,10).Value = 0.857
$ws.Cells.Item(20,11).Value = 0.8
$ws.Cells.Item(21,2).Value = 0.6909999999999999
$ws.Cells.Item(21,3).Value = 0.676
$ws.Cells.Item(21,4).Value = 0.787
$ws.Cells.Item(21,5).Value = 0.233
$ws.Cells.Item(21,6).Value = 0.507
$ws.Cells.Item(21,7).Value = 0.382
$ws.Cells.Item(21,8).Value = 0.436
$ws.Cells.Item(21,9).Value = 0.748
$ws.Cells.Item(21,10).Value = 0.832
$ws.Cells.Item(21,11).Value = 0.787
$ws.Cells.Item(22,2).Value = 0.674
$ws.Cells.Item(22,3).Value = 0.651
$ws.Cells.Item(22,4).Value = 0.778
$ws.Cells.Item(22,5).Value = 0.177
$ws.Cells.Item(22,6).Value = 0.468
$ws.Cells.Item(22,7).Value = 0.326
$ws.Cells.Item(22,8).Value = 0.384
$ws.Cells.Item(22,9).Value = 0.731
$ws.Cells.Item(22,10).Value = 0.832
$ws.Cells.Item(22,11).Value = 0.778
$ws.Cells.Item(23,2).Value = 0.6840000000000001
$ws.Cells.Item(23,3).Value = 0.651
$ws.Cells.Item(23,4).Value = 0.782
$ws.Cells.Item(23,5).Value = 0.22
$ws.Cells.Item(23,6).Value = 0.493
$ws.Cells.Item(23,7).Value = 0.382
$ws.Cells.Item(23,8).Value = 0.43
$ws.Cells.Item(23,9).Value = 0.745
$ws.Cells.Item(23,10).Value = 0.821
$ws.Cells.Item(23,11).Value = 0.782
$ws.Cells.Item(24,2).Value = 0.649
$ws.Cells.Item(24,3).Value = 0.6
$ws.Cells.Item(24,4).Value = 0.762
$ws.Cells.Item(24,5).Value = 0.11
$ws.Cells.Item(24,6).Value = 0.41
$ws.Cells.Item(24,7).Value = 0.281
$ws.Cells.Item(24,8).Value = 0.333
$ws.Cells.Item(24,9).Value = 0.714
$ws.Cells.Item(24,10).Value = 0.8159999999999999
$ws.Cells.Item(24,11).Value = 0.762
$ws.Cells.Item(25,2).Value = 0.667
$ws.Cells.Item(25,3).Value = 0.658
$ws.Cells.Item(25,4).Value = 0.764
$ws.Cells.Item(25,5).Value = 0.198
$ws.Cells.Item(25,6).Value = 0.462
$ws.Cells.Item(25,7).Value = 0.404
$ws.Cells.Item(25,8).Value = 0.431
$ws.Cells.Item(25,9).Value = 0.744
$ws.Cells.Item(25,10).Value = 0.786
$ws.Cells.Item(25,11).Value = 0.764

Write-Output "Updated 238 cells in B2:K25"